$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 59.42857
$ws.Range("I6").Value = 62.833332
$ws.Range("K6").Value = 188.499996
$ws.Range("M6").Value = -76.49999600000001
$ws.Range("H7").Value = 505
$ws.Range("I7").Value = 505
$ws.Range("K7").Value = 505
$ws.Range("M7").Value = -393
$ws.Range("H14").Value = 505
$ws.Range("I14").Value = 505
$ws.Range("K14").Value = 505
$ws.Range("M14").Value = -314
$ws.Range("H43").Value = 899
$ws.Range("J43").Value = 899
$ws.Range("L43").Value = 899
$ws.Range("N43").Value = -1037
$ws.Range("H62").Value = 3725
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3725
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H74").Value = 21617.555
$ws.Range("I74").Value = 36249.4
$ws.Range("J74").Value = 3327.75
$ws.Range("K74").Value = 36249.4
$ws.Range("L74").Value = 3327.75
$ws.Range("M74").Value = -35313.4
$ws.Range("N74").Value = -5199.75
$ws.Range("H77").Value = 21617.555
$ws.Range("I77").Value = 36249.4
$ws.Range("J77").Value = 3327.75
$ws.Range("K77").Value = 181247
$ws.Range("L77").Value = 16638.75
$ws.Range("M77").Value = -176567
$ws.Range("N77").Value = -25998.75
$ws.Range("H116").Value = 8477.4
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 8477.4
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 8477.4
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -15361.4
$ws.Range("H127").Value = 5581.75
$ws.Range("I127").Value = 5961
$ws.Range("K127").Value = 17883
$ws.Range("M127").Value = -12923
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 4688.1665
$ws.Range("I4").Value = 2377.3333
$ws.Range("J4").Value = 6999
$ws.Range("K4").Value = 2377.3333
$ws.Range("L4").Value = 6999
$ws.Range("M4").Value = -2261.3333
$ws.Range("N4").Value = -7231
$ws.Range("H110").Value = 2500
$ws.Range("I110").Value = 2500
$ws.Range("K110").Value = 2500
$ws.Range("M110").Value = -455
$ws.Range("H132").Value = 2471.25
$ws.Range("I132").Value = 2253.2632
$ws.Range("K132").Value = 6759.7896
$ws.Range("M132").Value = -4229.7896

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 39000
$ws.Range("J124").Value = 39000
$ws.Range("L124").Value = 39000
$ws.Range("N124").Value = -48820

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 70.67856999999999
$ws.Range("I7").Value = 47.266666
$ws.Range("J7").Value = 97.69231000000001
$ws.Range("K7").Value = 47.266666
$ws.Range("L7").Value = 97.69231000000001
$ws.Range("M7").Value = 65.733334
$ws.Range("N7").Value = -323.69231
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H122").Value = 1713
$ws.Range("I122").Value = 569.5
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 1708.5
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = 741.5
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 1537.9565
$ws.Range("I132").Value = 1413.65
$ws.Range("K132").Value = 4240.950000000001
$ws.Range("M132").Value = -1710.950000000001
$ws.Range("H133").Value = 37211.285
$ws.Range("I133").Value = 36080
$ws.Range("J133").Value = 43999
$ws.Range("K133").Value = 36080
$ws.Range("L133").Value = 43999
$ws.Range("M133").Value = -33550
$ws.Range("N133").Value = -49059

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7373
$ws.Range("I3").Value = 7373
$ws.Range("K3").Value = 22119
$ws.Range("M3").Value = -22007
$ws.Range("H4").Value = 7333794
$ws.Range("I4").Value = 7333794
$ws.Range("K4").Value = 22001382
$ws.Range("M4").Value = -22001270
$ws.Range("H12").Value = 18831.834
$ws.Range("J12").Value = 18831.834
$ws.Range("L12").Value = 56495.50199999999
$ws.Range("N12").Value = -56841.50199999999
$ws.Range("H25").Value = 1148.125
$ws.Range("I25").Value = 1230.8334
$ws.Range("J25").Value = 900
$ws.Range("K25").Value = 3692.5002
$ws.Range("L25").Value = 2700
$ws.Range("M25").Value = -3523.5002
$ws.Range("N25").Value = -3038
$ws.Range("H30").Value = 1148.125
$ws.Range("I30").Value = 1230.8334
$ws.Range("J30").Value = 900
$ws.Range("K30").Value = 3692.5002
$ws.Range("L30").Value = 2700
$ws.Range("M30").Value = -3590.5002
$ws.Range("N30").Value = -2904
$ws.Range("H42").Value = 6500
$ws.Range("J42").Value = 6500
$ws.Range("L42").Value = 19500
$ws.Range("N42").Value = -20568
$ws.Range("H60").Value = 973.5
$ws.Range("I60").Value = 982
$ws.Range("K60").Value = 2946
$ws.Range("M60").Value = -2695
$ws.Range("H92").Value = 126.27273
$ws.Range("J92").Value = 221.44444
$ws.Range("L92").Value = 664.33332
$ws.Range("N92").Value = -3160.33332
$ws.Range("H117").Value = 11576.777
$ws.Range("I117").Value = 399.16666
$ws.Range("J117").Value = 33932
$ws.Range("K117").Value = 1197.49998
$ws.Range("L117").Value = 101796
$ws.Range("M117").Value = 2244.50002
$ws.Range("N117").Value = -108680
$ws.Range("H123").Value = 7052.1113
$ws.Range("I123").Value = 3078.6667
$ws.Range("J123").Value = 14999
$ws.Range("K123").Value = 9236.000100000001
$ws.Range("L123").Value = 44997
$ws.Range("M123").Value = -6786.000100000001
$ws.Range("N123").Value = -49897

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 688241.9
$ws.Range("I3").Value = 917006
$ws.Range("J3").Value = 1949.5
$ws.Range("K3").Value = 917006
$ws.Range("L3").Value = 1949.5
$ws.Range("M3").Value = -916890
$ws.Range("N3").Value = -2181.5
$ws.Range("H11").Value = 898649.75
$ws.Range("I11").Value = 1104170.8
$ws.Range("K11").Value = 1104170.8
$ws.Range("M11").Value = -1104031.8
$ws.Range("H14").Value = 1175.875
$ws.Range("I14").Value = 401.16666
$ws.Range("K14").Value = 401.16666
$ws.Range("M14").Value = -233.16666
$ws.Range("H70").Value = 9833.333000000001
$ws.Range("I70").Value = 9000
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 10000
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -10540
$ws.Range("H73").Value = 9833.333000000001
$ws.Range("I73").Value = 9000
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 10000
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -11872
$ws.Range("H80").Value = 1800
$ws.Range("I80").Value = 1700
$ws.Range("J80").Value = 1900
$ws.Range("K80").Value = 1700
$ws.Range("L80").Value = 1900
$ws.Range("M80").Value = -702
$ws.Range("N80").Value = -3896
$ws.Range("H83").Value = 1800
$ws.Range("I83").Value = 1700
$ws.Range("J83").Value = 1900
$ws.Range("K83").Value = 8500
$ws.Range("L83").Value = 9500
$ws.Range("M83").Value = -3508
$ws.Range("N83").Value = -19484
$ws.Range("H135").Value = 526315
$ws.Range("J135").Value = 526315
$ws.Range("L135").Value = 526315
$ws.Range("N135").Value = -536455

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H61").Value = 3527.7
$ws.Range("I61").Value = 3988
$ws.Range("K61").Value = 3988
$ws.Range("M61").Value = -3786
$ws.Range("H82").Value = 1202.6666
$ws.Range("I82").Value = 1202.6666
$ws.Range("K82").Value = 1202.6666
$ws.Range("M82").Value = -841.6666
$ws.Range("H85").Value = 1202.6666
$ws.Range("I85").Value = 1202.6666
$ws.Range("K85").Value = 1202.6666
$ws.Range("M85").Value = 45.33339999999998
$ws.Range("H104").Value = 22872
$ws.Range("J104").Value = 22872
$ws.Range("L104").Value = 22872
$ws.Range("N104").Value = -29860
$ws.Range("H113").Value = 3527.7
$ws.Range("I113").Value = 3988
$ws.Range("K113").Value = 3988
$ws.Range("M113").Value = -1818
$ws.Range("H122").Value = 7544.8184
$ws.Range("I122").Value = 7799.467
$ws.Range("K122").Value = 23398.401
$ws.Range("M122").Value = -20948.401

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 300
$ws.Range("I6").Value = 300
$ws.Range("K6").Value = 300
$ws.Range("M6").Value = -185
$ws.Range("H14").Value = 1000
$ws.Range("J14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1336

